# Update vm_pu results for Case_1_35 (380 kV slack-bus case): res_bus/vm_pu
# Column B (slack/Vm setpoint) drops from 1.05 p.u. to 1.02 p.u.; all other
# bus voltage magnitudes are refreshed with the corresponding re-run output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (bus index 0)
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037092216978502
$ws.Range("D2").Value = 1.060815371643894
$ws.Range("E2").Value = 1.048074343306386
$ws.Range("F2").Value = 1.064150523653378
$ws.Range("I2").Value = 1.04895077201893
$ws.Range("J2").Value = 1.042196970815508
$ws.Range("K2").Value = 1.063541124872434
$ws.Range("L2").Value = 1.050835203594305
$ws.Range("M2").Value = 1.06686723736307
$ws.Range("N2").Value = 1.018010165150772

# Row 3 (bus index 1)
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037887722550557
$ws.Range("D3").Value = 1.06145531354983
$ws.Range("E3").Value = 1.048767079419104
$ws.Range("F3").Value = 1.064891529607091
$ws.Range("I3").Value = 1.049164905398065
$ws.Range("J3").Value = 1.042637486616989
$ws.Range("K3").Value = 1.063996006584921
$ws.Range("L3").Value = 1.051340220795521
$ws.Range("M3").Value = 1.067423578369567
$ws.Range("N3").Value = 1.018157050878827

# Row 4 (bus index 2)
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038403187605517
$ws.Range("D4").Value = 1.061869757957021
$ws.Range("E4").Value = 1.049216286123253
$ws.Range("F4").Value = 1.065371853381174
$ws.Range("I4").Value = 1.04930244466281
$ws.Range("J4").Value = 1.042922562915059
$ws.Range("K4").Value = 1.064290014229189
$ws.Range("L4").Value = 1.051667281315053
$ws.Range("M4").Value = 1.067783740624427
$ws.Range("N4").Value = 1.018252080372798

# Row 5 (bus index 3)
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038620059588229
$ws.Range("D5").Value = 1.06204407430002
$ws.Range("E5").Value = 1.049405360481597
$ws.Range("F5").Value = 1.065573981239008
$ws.Range("I5").Value = 1.049360021178175
$ws.Range("J5").Value = 1.043042415595882
$ws.Range("K5").Value = 1.064413534133785
$ws.Range("L5").Value = 1.051804843210279
$ws.Range("M5").Value = 1.067935192361303
$ws.Range("N5").Value = 1.018292026609366

# Row 6 (bus index 4)
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038656483291677
$ws.Range("D6").Value = 1.062073347641271
$ws.Range("E6").Value = 1.049437120202105
$ws.Range("F6").Value = 1.065607931045865
$ws.Range("I6").Value = 1.04936967413422
$ws.Range("J6").Value = 1.043062539737203
$ws.Range("K6").Value = 1.064434268877906
$ws.Range("L6").Value = 1.05182794426323
$ws.Range("M6").Value = 1.067960624041087
$ws.Range("N6").Value = 1.018298733500266

# Row 7 (bus index 5)
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03840608479202
$ws.Range("D7").Value = 1.061872086853421
$ws.Range("E7").Value = 1.04921881165126
$ws.Range("F7").Value = 1.065374553441369
$ws.Range("I7").Value = 1.049303214966873
$ws.Range("J7").Value = 1.042924164368154
$ws.Range("K7").Value = 1.064291665028139
$ws.Range("L7").Value = 1.05166911916778
$ws.Range("M7").Value = 1.067785764176846
$ws.Range("N7").Value = 1.018252614153322

# Row 8 (bus index 6)
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037360911675677
$ws.Range("D8").Value = 1.061031567382713
$ws.Range("E8").Value = 1.048308256722353
$ws.Range("F8").Value = 1.064400774160046
$ws.Range("I8").Value = 1.049023349884625
$ws.Range("J8").Value = 1.042345837540323
$ws.Range("K8").Value = 1.063694921895433
$ws.Range("L8").Value = 1.051005817499101
$ws.Range("M8").Value = 1.067055218656272
$ws.Range("N8").Value = 1.018059808770892

# Row 9 (bus index 7)
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035524776916731
$ws.Range("D9").Value = 1.059553307622953
$ws.Range("E9").Value = 1.046711184082629
$ws.Range("F9").Value = 1.06269140542562
$ws.Range("I9").Value = 1.048522427193618
$ws.Range("J9").Value = 1.041327063939484
$ws.Range("K9").Value = 1.062640918470752
$ws.Range("L9").Value = 1.049839211990815
$ws.Range("M9").Value = 1.065769302154028
$ws.Range("N9").Value = 1.017719964248749

# Row 10 (bus index 8)
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034304557811647
$ws.Range("D10").Value = 1.05856983673471
$ws.Range("E10").Value = 1.045651591387235
$ws.Range("F10").Value = 1.061556362500615
$ws.Range("I10").Value = 1.048183315394869
$ws.Range("J10").Value = 1.040648170448369
$ws.Range("K10").Value = 1.061936687766614
$ws.Range("L10").Value = 1.049063056716841
$ws.Range("M10").Value = 1.06491306917612
$ws.Range("N10").Value = 1.017493364644664

# Row 11 (bus index 9)
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033777130742442
$ws.Range("D11").Value = 1.058144492108757
$ws.Range("E11").Value = 1.045194016275675
$ws.Range("F11").Value = 1.061065979685741
$ws.Range("I11").Value = 1.048035263640485
$ws.Range("J11").Value = 1.040354287310121
$ws.Range("K11").Value = 1.061631398000886
$ws.Range("L11").Value = 1.048727367407224
$ws.Range("M11").Value = 1.064542580800266
$ws.Range("N11").Value = 1.017395241864037

# Row 12 (bus index 10)
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033581363321515
$ws.Range("D12").Value = 1.057986578147261
$ws.Range("E12").Value = 1.045024240188311
$ws.Range("F12").Value = 1.060883996903765
$ws.Range("I12").Value = 1.047980089181951
$ws.Range("J12").Value = 1.040245139731113
$ws.Range("K12").Value = 1.061517948362093
$ws.Range("L12").Value = 1.048602737769053
$ws.Range("M12").Value = 1.064405006704927
$ws.Range("N12").Value = 1.017358794632384

# Row 13 (bus index 11)
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033623349620386
$ws.Range("D13").Value = 1.058020447665434
$ws.Range("E13").Value = 1.045060649209299
$ws.Range("F13").Value = 1.060923025225014
$ws.Range("I13").Value = 1.04799193248593
$ws.Range("J13").Value = 1.040268551615525
$ws.Range("K13").Value = 1.061542285995571
$ws.Range("L13").Value = 1.048629468494839
$ws.Range("M13").Value = 1.064434514895769
$ws.Range("N13").Value = 1.017366612682284

# Row 14 (bus index 12)
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033760945629038
$ws.Range("D14").Value = 1.05813143728851
$ws.Range("E14").Value = 1.045179978690764
$ws.Range("F14").Value = 1.061050933507903
$ws.Range("I14").Value = 1.048030706598392
$ws.Range("J14").Value = 1.040345264849833
$ws.Range("K14").Value = 1.06162202125796
$ws.Range("L14").Value = 1.048717064242926
$ws.Range("M14").Value = 1.064531208019877
$ws.Range("N14").Value = 1.017392229122948

# Row 15 (bus index 13)
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033845742006801
$ws.Range("D15").Value = 1.05819983205417
$ws.Range("E15").Value = 1.045253526453985
$ws.Range("F15").Value = 1.061129764248654
$ws.Range("I15").Value = 1.048054572594329
$ws.Range("J15").Value = 1.040392532261599
$ws.Range("K15").Value = 1.06167114201922
$ws.Range("L15").Value = 1.048771042913926
$ws.Range("M15").Value = 1.064590789435332
$ws.Range("N15").Value = 1.01740801226384

# Row 16 (bus index 14)
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034339581286561
$ws.Range("D16").Value = 1.058598076273868
$ws.Range("E16").Value = 1.045681985347268
$ws.Range("F16").Value = 1.061588930922407
$ws.Range("I16").Value = 1.04819311560437
$ws.Range("J16").Value = 1.040667676352029
$ws.Range("K16").Value = 1.061956941541029
$ws.Range("L16").Value = 1.04908534368351
$ws.Range("M16").Value = 1.064937663030438
$ws.Range("N16").Value = 1.017499876695028

# Row 17 (bus index 15)
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034649605497526
$ws.Range("D17").Value = 1.058848020999317
$ws.Range("E17").Value = 1.045951078529841
$ws.Range("F17").Value = 1.061877249655842
$ws.Range("I17").Value = 1.048279695601362
$ws.Range("J17").Value = 1.040840289937949
$ws.Range("K17").Value = 1.062136122621439
$ws.Range("L17").Value = 1.049282602053972
$ws.Range("M17").Value = 1.065155320198474
$ws.Range("N17").Value = 1.017557500190922

# Row 18 (bus index 16)
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034830527545267
$ws.Range("D18").Value = 1.058993858087746
$ws.Range("E18").Value = 1.046108155102136
$ws.Range("F18").Value = 1.062045526912173
$ws.Range("I18").Value = 1.048330079040511
$ws.Range("J18").Value = 1.04094098036507
$ws.Range("K18").Value = 1.062240601723229
$ws.Range("L18").Value = 1.049397697149335
$ws.Range("M18").Value = 1.065282301501732
$ws.Range("N18").Value = 1.017591110615693

# Row 19 (bus index 17)
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034892232525608
$ws.Range("D19").Value = 1.059043592937482
$ws.Range("E19").Value = 1.046161734287948
$ws.Range("F19").Value = 1.062102922988478
$ws.Range("I19").Value = 1.048347238581452
$ws.Range("J19").Value = 1.040975314487826
$ws.Range("K19").Value = 1.062276220547631
$ws.Range("L19").Value = 1.049436947936695
$ws.Range("M19").Value = 1.065325603125306
$ws.Range("N19").Value = 1.01760257082175

# Row 20 (bus index 18)
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034616333493434
$ws.Range("D20").Value = 1.058821199257651
$ws.Range("E20").Value = 1.045922195047032
$ws.Range("F20").Value = 1.061846304823135
$ws.Range("I20").Value = 1.048270418503542
$ws.Range("J20").Value = 1.040821769312815
$ws.Range("K20").Value = 1.062116901711952
$ws.Range("L20").Value = 1.049261434194595
$ws.Range("M20").Value = 1.065131964987519
$ws.Range("N20").Value = 1.017551317771511

# Row 21 (bus index 19)
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033720423085875
$ws.Range("D21").Value = 1.058098751441335
$ws.Range("E21").Value = 1.045144833935197
$ws.Range("F21").Value = 1.061013263072384
$ws.Range("I21").Value = 1.048019293589739
$ws.Range("J21").Value = 1.040322674325226
$ws.Range("K21").Value = 1.061598542616996
$ws.Range("L21").Value = 1.048691267800728
$ws.Range("M21").Value = 1.064502733136993
$ws.Range("N21").Value = 1.017384685717447

# Row 22 (bus index 20)
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033157953907037
$ws.Range("D22").Value = 1.057644971996528
$ws.Range("E22").Value = 1.044657163041529
$ws.Range("F22").Value = 1.060490465300668
$ws.Range("I22").Value = 1.047860352081028
$ws.Range("J22").Value = 1.040008953450693
$ws.Range("K22").Value = 1.061272332951849
$ws.Range("L22").Value = 1.048333131314974
$ws.Range("M22").Value = 1.06410735273233
$ws.Range("N22").Value = 1.017279917413461

# Row 23 (bus index 21)
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033456050610137
$ws.Range("D23").Value = 1.057885485492986
$ws.Range("E23").Value = 1.044915582897701
$ws.Range("F23").Value = 1.060767517732051
$ws.Range("I23").Value = 1.047944709123275
$ws.Range("J23").Value = 1.040175254792699
$ws.Range("K23").Value = 1.061445290446403
$ws.Range("L23").Value = 1.048522952566429
$ws.Range("M23").Value = 1.064316927767412
$ws.Range("N23").Value = 1.017335456922295

# Row 24 (bus index 22)
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034631367401535
$ws.Range("D24").Value = 1.058833318698427
$ws.Range("E24").Value = 1.045935245882005
$ws.Range("F24").Value = 1.061860287134371
$ws.Range("I24").Value = 1.048274610786995
$ws.Range("J24").Value = 1.040830137962045
$ws.Range("K24").Value = 1.062125586919453
$ws.Range("L24").Value = 1.049270998923039
$ws.Range("M24").Value = 1.065142518123299
$ws.Range("N24").Value = 1.017554111341979

# Row 25 (bus index 23)
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035998787564307
$ws.Range("D25").Value = 1.059935123491437
$ws.Range("E25").Value = 1.047123171496796
$ws.Range("F25").Value = 1.063132528120746
$ws.Range("I25").Value = 1.048652841979254
$ws.Range("J25").Value = 1.041590396196421
$ws.Range("K25").Value = 1.062913685919552
$ws.Range("L25").Value = 1.050140535364892
$ws.Range("M25").Value = 1.06610156591913
$ws.Range("N25").Value = 1.017807830737182
